$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.Formula = $val
    $cell.Copy()
    $cell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

Set-TextValue 'D2' '="27.502.20"'
Set-TextValue 'E2' '="  -0.16%  "'
Set-TextValue 'D3' '="1.580.57"'
Set-TextValue 'E3' '="  -0.91%  "'
Set-TextValue 'E4' '="  +0.03%  "'
Set-TextValue 'D5' '="207.58"'
Set-TextValue 'E5' '="  -0.13%  "'
Set-TextValue 'E6' '="  -0.98%  "'
Set-TextValue 'E7' '="  +0.06%  "'
Set-TextValue 'D8' '="22.20"'
Set-TextValue 'E8' '="  -0.25%  "'
Set-TextValue 'E9' '="  -0.88%  "'
Set-TextValue 'D10' '="0.0589"'
Set-TextValue 'D11' '="0.0865"'
Set-TextValue 'E11' '="  -0.58%  "'
Set-TextValue 'D12' '="1.804.28"'
Set-TextValue 'E12' '="  -0.89%  "'
Set-TextValue 'D13' '="1.589.93"'
Set-TextValue 'E13' '="  -0.58%  "'
Set-TextValue 'E14' '="  -1.22%  "'
Set-TextValue 'E15' '="  -2.82%  "'
Set-TextValue 'D16' '="27.492.63"'
Set-TextValue 'E16' '="  -0.20%  "'
Set-TextValue 'D17' '="62.94"'
Set-TextValue 'E17' '="  -0.64%  "'
Set-TextValue 'D18' '="215.03"'
Set-TextValue 'E18' '="  -1.78%  "'
Set-TextValue 'D19' '="0.0₃0690"'
Set-TextValue 'E19' '="  -0.62%  "'
Set-TextValue 'D20' '="7.28"'
Set-TextValue 'E20' '="  -1.42%  "'
Set-TextValue 'E21' '="  +0.02%  "'
Set-TextValue 'D22' '="4.13"'
Set-TextValue 'E22' '="  -1.91%  "'
Set-TextValue 'D23' '="9.70"'
Set-TextValue 'E23' '="  +0.52%  "'
Set-TextValue 'E24' '="  +0.38%  "'
Set-TextValue 'D25' '="152.99"'
Set-TextValue 'E25' '="  -1.08%  "'
Set-TextValue 'E26' '="  +2.46%  "'
Set-TextValue 'E27' '="  +0.05%  "'
Set-TextValue 'D28' '="15.04"'
Set-TextValue 'E28' '="  -0.15%  "'
Set-TextValue 'E29' '="  -1.25%  "'
Set-TextValue 'E30' '="  -0.30%  "'
Set-TextValue 'E31' '="  +0.52%  "'
Set-TextValue 'D32' '="3.22"'
Set-TextValue 'E32' '="  -1.48%  "'
Set-TextValue 'D33' '="1.364.99"'
Set-TextValue 'E33' '="  +0.23%  "'
Set-TextValue 'D34' '="2.95"'
Set-TextValue 'E34' '="  -0.30%  "'
Set-TextValue 'E35' '="  -0.45%  "'
Set-TextValue 'D36' '="0.971"'
Set-TextValue 'E36' '="  +0.67%  "'
Set-TextValue 'E37' '="  +0.30%  "'
Set-TextValue 'E38' '="  +1.10%  "'
Set-TextValue 'D39' '="0.530"'
Set-TextValue 'E39' '="  -1.67%  "'
Set-TextValue 'D40' '="0.821"'
Set-TextValue 'E40' '="  +0.93%  "'
Set-TextValue 'E41' '="  +0.06%  "'
Set-TextValue 'D42' '="0.971"'
Set-TextValue 'E42' '="  +0.01%  "'
Set-TextValue 'D43' '="64.09"'
Set-TextValue 'E43' '="  +0.11%  "'
Set-TextValue 'E44' '="  +3.62%  "'
Set-TextValue 'E45' '="  -2.08%  "'
Set-TextValue 'E46' '="  -1.13%  "'
Set-TextValue 'D47' '="1.716.35"'
Set-TextValue 'E47' '="  -0.76%  "'
Set-TextValue 'D48' '="86.35"'
Set-TextValue 'E48' '="  -1.60%  "'
Set-TextValue 'D49' '="0.0₇0999"'
Set-TextValue 'E49' '="  -0.23%  "'
Set-TextValue 'D50' '="0.0956"'
Set-TextValue 'E50' '="  -1.75%  "'
Set-TextValue 'E51' '="  -0.80%  "'
